$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Day 2 (20250701) data for plate15, wells A02:D06 (20 rows), appended
# immediately after the existing data (which ends at row 461).
$newData = @(
    @(20250701, "plate15", "A02", 16.59),
    @(20250701, "plate15", "A03", 14.06),
    @(20250701, "plate15", "A04", 13.44),
    @(20250701, "plate15", "A05", 14.97),
    @(20250701, "plate15", "A06", 16.36),
    @(20250701, "plate15", "B02", 15.82),
    @(20250701, "plate15", "B03", 12.16),
    @(20250701, "plate15", "B04", 14.12),
    @(20250701, "plate15", "B05", 13.39),
    @(20250701, "plate15", "B06", 17.5),
    @(20250701, "plate15", "C02", 12.9),
    @(20250701, "plate15", "C03", 14.49),
    @(20250701, "plate15", "C04", 14.98),
    @(20250701, "plate15", "C05", 14.72),
    @(20250701, "plate15", "C06", 13.38),
    @(20250701, "plate15", "D02", 13.68),
    @(20250701, "plate15", "D03", 16.63),
    @(20250701, "plate15", "D04", 12.58),
    @(20250701, "plate15", "D05", 13.47),
    @(20250701, "plate15", "D06", 13.81)
)

# Find the first empty row after the existing data (last used row is 461).
$startRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1

for ($i = 0; $i -lt $newData.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $newData[$i]
    $ws.Cells.Item($r, 1).Value = $rowVals[0]
    $ws.Cells.Item($r, 2).Value = $rowVals[1]
    $ws.Cells.Item($r, 3).Value = $rowVals[2]
    $ws.Cells.Item($r, 4).Value = $rowVals[3]
}

$lastRow = $startRow + $newData.Count - 1

# Match the saved view state: scrolled down near the bottom of the sheet,
# with G469 as the active/selected cell.
$app = $excel
$topLeftRow = $lastRow - 12
$app.Goto($ws.Cells.Item($topLeftRow, 1), $true)
$ws.Range("G469").Select()
